# lab3 final final fix
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in header: "квартиросъемщика" -> "квартиросъёмщика"
$ws.Range("B2").Value = "Фамилия квартиросъёмщика"

# Rename tenant entry "Куропаткин 1" -> "Баницин"
$ws.Range("B37").Value = "Баницин"

# Shorten label for total sum row
$ws.Range("B40").Value = "Общая сумма, руб."

# Widen column B to fit new content (~31.125 chars)
$ws.Columns.Item(2).ColumnWidth = 30.43

# Reset row 40 height to default (remove custom 30.75 height / autosize it back)
$ws.Rows.Item(40).AutoFit()

# Reset view: scroll back to top-left and select B1 (previously was topLeftCell A31 / selection B45)
$ws.Range("B1").Select()
